$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "dnasr281@gmail.com, *") {
        $rest = $val.Substring(20)
        $cell.Value2 = "$rest, dnasr281@gmail.com"
    }
}
